$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5788
$ws.Range("B2").Value = 45828
$ws.Range("A3").Value = 5719
$ws.Range("B3").Value = 45828.01041666666
$ws.Range("A4").Value = 5682
$ws.Range("B4").Value = 45828.02083333334
$ws.Range("A5").Value = 5639
$ws.Range("B5").Value = 45828.03125
$ws.Range("A6").Value = 5610
$ws.Range("B6").Value = 45828.04166666666
$ws.Range("A7").Value = 5518
$ws.Range("B7").Value = 45828.05208333334
$ws.Range("A8").Value = 5500
$ws.Range("B8").Value = 45828.0625
$ws.Range("A9").Value = 5513
$ws.Range("B9").Value = 45828.07291666666
$ws.Range("A10").Value = 5495
$ws.Range("B10").Value = 45828.08333333334
$ws.Range("A11").Value = 5513
$ws.Range("B11").Value = 45828.09375
$ws.Range("A12").Value = 5453
$ws.Range("B12").Value = 45828.10416666666
$ws.Range("A13").Value = 5549
$ws.Range("B13").Value = 45828.11458333334
$ws.Range("A14").Value = 5509
$ws.Range("B14").Value = 45828.125
$ws.Range("A15").Value = 5475
$ws.Range("B15").Value = 45828.13541666666
$ws.Range("A16").Value = 5454
$ws.Range("B16").Value = 45828.14583333334
$ws.Range("A17").Value = 5473
$ws.Range("B17").Value = 45828.15625
$ws.Range("A18").Value = 5402
$ws.Range("B18").Value = 45828.16666666666
$ws.Range("A19").Value = 5376
$ws.Range("B19").Value = 45828.17708333334
$ws.Range("A20").Value = 5391
$ws.Range("B20").Value = 45828.1875
$ws.Range("A21").Value = 5443
$ws.Range("B21").Value = 45828.19791666666
$ws.Range("A22").Value = 5524
$ws.Range("B22").Value = 45828.20833333334
$ws.Range("A23").Value = 5703
$ws.Range("B23").Value = 45828.21875
$ws.Range("A24").Value = 5734
$ws.Range("B24").Value = 45828.22916666666
$ws.Range("A25").Value = 5779
$ws.Range("B25").Value = 45828.23958333334
$ws.Range("A26").Value = 5984
$ws.Range("B26").Value = 45828.25
$ws.Range("A27").Value = 6081
$ws.Range("B27").Value = 45828.26041666666
$ws.Range("A28").Value = 6163
$ws.Range("B28").Value = 45828.27083333334
$ws.Range("A29").Value = 6231
$ws.Range("B29").Value = 45828.28125
$ws.Range("A30").Value = 6285
$ws.Range("B30").Value = 45828.29166666666
$ws.Range("A31").Value = 6271
$ws.Range("B31").Value = 45828.30208333334
$ws.Range("A32").Value = 6225
$ws.Range("B32").Value = 45828.3125
$ws.Range("A33").Value = 6209
$ws.Range("B33").Value = 45828.32291666666
$ws.Range("A34").Value = 6062
$ws.Range("B34").Value = 45828.33333333334
$ws.Range("A35").Value = 5962
$ws.Range("B35").Value = 45828.34375
$ws.Range("A36").Value = 5975
$ws.Range("B36").Value = 45828.35416666666
$ws.Range("A37").Value = 5860
$ws.Range("B37").Value = 45828.36458333334
$ws.Range("A38").Value = 5671
$ws.Range("B38").Value = 45828.375
$ws.Range("A39").Value = 5639
$ws.Range("B39").Value = 45828.38541666666
$ws.Range("A40").Value = 5646
$ws.Range("B40").Value = 45828.39583333334
$ws.Range("A41").Value = 5651
$ws.Range("B41").Value = 45828.40625
$ws.Range("A42").Value = 5404
$ws.Range("B42").Value = 45828.41666666666
$ws.Range("A43").Value = 5303
$ws.Range("B43").Value = 45828.42708333334
$ws.Range("A44").Value = 5334
$ws.Range("B44").Value = 45828.4375
$ws.Range("A45").Value = 5307
$ws.Range("B45").Value = 45828.44791666666
$ws.Range("A46").Value = 5293
$ws.Range("B46").Value = 45828.45833333334
$ws.Range("A47").Value = 5225
$ws.Range("B47").Value = 45828.46875

# Apply the date number format to the newly appended rows (42-47)
# so their style matches the rest of column B (style index 2).
$ws.Range("B42:B47").NumberFormat = "YYYY-MM-DD HH:MM:SS"

